$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "R 1 300 000"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "74 m²"
$ws.Range("E2").Value = "Sherwood"

$ws.Range("A3").Value = "R 1 350 000"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "111 m²"
$ws.Range("E3").Value = "Essenwood"

$ws.Range("A4").Value = "R 1 495 000"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "90 m²"
$ws.Range("E4").Value = "North Beach"

$ws.Range("A5").Value = "R 1 500 000"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "77 m²"
$ws.Range("E5").Value = "North Beach"

$ws.Range("A6").Value = "R 1 500 000"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "77 m²"
$ws.Range("E6").Value = "North Beach"

$ws.Range("A7").Value = "R 1 690 000"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "146 m²"
$ws.Range("E7").Value = "Morningside"

$ws.Range("A8").Value = "R 1 700 000"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "164 m²"
$ws.Range("E8").Value = "Glenwood"

$ws.Range("A9").Value = "R 2 200 000"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "210 m²"
$ws.Range("E9").Value = "South Beach"

$ws.Range("A10").Value = "R 2 325 000"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "296 m²"
$ws.Range("E10").Value = "Glenwood"

$ws.Range("A11").Value = "R 3 395 000"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 2.5
$ws.Range("D11").Value = "209 m²"
$ws.Range("E11").Value = "Essenwood"

$ws.Range("A12").Value = "R 499 000"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "73 m²"
$ws.Range("E12").Value = "Umbilo"

$ws.Range("A13").Value = "R 550 000"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "51 m²"
$ws.Range("E13").Value = "Esplanade"

$ws.Range("A14").Value = "R 595 000"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "46 m²"
$ws.Range("E14").Value = "Morningside"

$ws.Range("A15").Value = "R 650 000"
$ws.Range("B15").Value = 2.5
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "81 m²"
$ws.Range("E15").Value = "Montclair"

$ws.Range("A16").Value = "R 725 000"
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = "565 m²"
$ws.Range("E16").Value = "Umbilo"

$ws.Range("A17").Value = "R 920 000"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "71 m²"
$ws.Range("E17").Value = "Musgrave"
